$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Names-key" (sheet1): just the cursor moves from D7 to H1 as part of
# this edit (the tab-selected flag will move to ValueSets once we activate
# it below).
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Names-key")
$null = $ws1.Range("H1").Select()

# ---------------------------------------------------------------------------
# Sheet "ValueSets" (sheet4): populate the new value-set tracker.
# Columns A/B mirror the "No" / "Core Profile Title" columns from the
# Names-key sheet (one row per core profile), columns C/D are new and list
# the value sets used per profile (only AllergyIntolerance has any so far -
# the 5 substance value sets), and E18 captures a build-tool error note.
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("ValueSets")

# Column widths (approximate the widened columns used for the new data)
$ws4.Columns.Item(1).ColumnWidth = 12.02
$ws4.Columns.Item(2).ColumnWidth = 31.88
$ws4.Columns.Item(3).ColumnWidth = 28.17
$ws4.Columns.Item(4).ColumnWidth = 98.59
$ws4.Columns.Item(5).ColumnWidth = 21.74

# Header row
$ws4.Range("A1").Value = "No"
$ws4.Range("A1").Font.Name = "Calibri"
$ws4.Range("A1").Font.Bold = $true
$ws4.Range("A1").Font.Size = 12
$ws4.Range("A1").VerticalAlignment = -4108
$ws4.Range("A1").Borders.Item(9).LineStyle = 1
$ws4.Range("A1").Borders.Item(9).Weight = -4138

$ws4.Range("B1").Value = "Core Profile Title "
$ws4.Range("B1").Font.Bold = $true
$ws4.Range("B1").WrapText = $true

$ws4.Range("C1").Value = "Valueset id"
$ws4.Range("C1").Font.Bold = $true

$ws4.Range("D1").Value = "Title"
$ws4.Range("D1").Font.Bold = $true

# Data rows: A = No., B = Core Profile Title (copied from Names-key)
$ws4.Range("A2").Value = 1
$ws4.Range("B2").Value = "AllergyIntolerance"
$ws4.Range("C2").Value = "daf-core-substance-ndfrt"
$ws4.Range("D2").Value = "DAF Core Substance ND-FRT codes"

$ws4.Range("A3").Value = 1
$ws4.Range("B3").Value = "AllergyIntolerance"
$ws4.Range("C3").Value = "daf-core-substance-rxnorm"
$ws4.Range("D3").Value = "DAF Core Substance RxNorm Codes"

$ws4.Range("A4").Value = 1
$ws4.Range("B4").Value = "AllergyIntolerance"
$ws4.Range("C4").Value = "daf-core-substance-sct"
$ws4.Range("D4").Value = "DAF Core SNOMED CT Substances Other Than Clinical Drugs"

$ws4.Range("A5").Value = 1
$ws4.Range("B5").Value = "AllergyIntolerance"
$ws4.Range("C5").Value = "daf-core-substance-unii"
$ws4.Range("D5").Value = "DAF Core Substance UNII Codes"

$ws4.Range("A6").Value = 1
$ws4.Range("B6").Value = "AllergyIntolerance"
$ws4.Range("C6").Value = "daf-core-substance"
$ws4.Range("D6").Value = "DAF Core Substance-Reactant for Intolerance and Negation Codes"

$ws4.Range("A7").Value = 2
$ws4.Range("B7").Value = "CarePlan"

$ws4.Range("A8").Value = 3
$ws4.Range("B8").Value = "CareTeam"

$ws4.Range("A9").Value = 4
$ws4.Range("B9").Value = "Conformance"

$ws4.Range("A10").Value = 5
$ws4.Range("B10").Value = "DocumentReference"

$ws4.Range("A11").Value = 6
$ws4.Range("B11").Value = "?"

$ws4.Range("A12").Value = 7
$ws4.Range("B12").Value = "Goals"

$ws4.Range("A13").Value = 8
$ws4.Range("B13").Value = "Immunization"

$ws4.Range("A14").Value = 9
$ws4.Range("B14").Value = "Device-UDI"

$ws4.Range("A15").Value = 10
$ws4.Range("B15").Value = "'-"

$ws4.Range("A16").Value = 11
$ws4.Range("B16").Value = "DiagnosticReport-Results"

$ws4.Range("A17").Value = 12
$ws4.Range("B17").Value = "Observation-Results"

$ws4.Range("A18").Value = 13
$ws4.Range("B18").Value = "Location"
$ws4.Range("E18").Value = "XPDY0002: Cannot evaluate function 'fn:root', because the context node is undefined"

$ws4.Range("A19").Value = 14
$ws4.Range("B19").Value = "'-"

$ws4.Range("A20").Value = 15
$ws4.Range("B20").Value = "MedicationOrder"

$ws4.Range("A21").Value = 16
$ws4.Range("B21").Value = "MedicationStatement"

$ws4.Range("A22").Value = 17
$ws4.Range("B22").Value = "Medication"

$ws4.Range("A23").Value = 18
$ws4.Range("B23").Value = "Organization"

$ws4.Range("A24").Value = 19
$ws4.Range("B24").Value = "Patient"

$ws4.Range("A25").Value = 20
$ws4.Range("B25").Value = "Practitioner"

$ws4.Range("A26").Value = 21
$ws4.Range("B26").Value = "Condition"

$ws4.Range("A27").Value = 22
$ws4.Range("B27").Value = "Procedure"

$ws4.Range("A28").Value = 23
$ws4.Range("B28").Value = "Observation-Smokingstatus"

$ws4.Range("A29").Value = 24
$ws4.Range("B29").Value = "Observation-Vitalsigns"

$ws4.Range("B32").Value = "MedicationDispense"
$ws4.Range("B33").Value = "MedicationAdministration"
$ws4.Range("B34").Value = "Observation-Resultsv2"

# Bold + centered "No." column (matches the style used on Names-key col A)
$ws4.Range("A2:A29").Font.Bold = $true
$ws4.Range("A2:A29").HorizontalAlignment = -4108

# ValueSets becomes the active / selected tab, cursor parked at D14
$ws4.Activate()
$null = $ws4.Range("D14").Select()
